# Adds two new "rule" paragraphs after the "Adjectives ... etc." paragraph,
# and relocates the hidden "_GoBack" bookmark (which Word keeps pinned to the
# most-recently-edited spot) to the end of the very last paragraph we insert -
# matching what real Word does when you type new paragraphs at that location.

$d = $word.ActiveDocument

# --- locate the "Adjectives - old, interesting, expensive, etc." paragraph ---
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Adjectives*old, interesting, expensive, etc.*") {
        $targetIndex = $i
    }
}

$p = $d.Paragraphs($targetIndex)
$r = $p.Range
$r.InsertParagraphAfter()

$p2 = $d.Paragraphs($targetIndex + 1)
$p2.Range.InsertBefore("Basic word order in English")

$p2 = $d.Paragraphs($targetIndex + 1)
$r2 = $p2.Range
$r2.InsertParagraphAfter()

$p3 = $d.Paragraphs($targetIndex + 2)
$p3.Range.InsertBefore("Adverbs of frequency with present simple")

# --- move the "_GoBack" bookmark to sit right after the text we just typed ---
# (collapsed ranges that land exactly on "paragraph end - 1" confuse
# Bookmarks.Add in this host, so nudge the boundary out of the way with a
# throwaway run, place the bookmark at the now mid-paragraph offset, then
# delete the throwaway text again.)
$pLast = $d.Paragraphs($targetIndex + 2)
$prLast = $pLast.Range
$origEnd = $prLast.End - 1
$prLast.InsertAfter("ZZTEMPZZ")

$rBm = $d.Range($origEnd, $origEnd)
$d.Bookmarks.Add("_GoBack", $rBm)

$delRange = $d.Range($origEnd, $origEnd + 8)
$delRange.Delete()
